$d = $word.ActiveDocument

# The document repeats the same interviewer placeholder block (Họ & tên /
# Công việc / Tuổi / Giới tính) once per interviewer. The diff only touches
# the block that belongs to "Người phỏng vấn thứ 5" (the 5th interviewer),
# so first locate that heading paragraph, then edit the four paragraphs
# that immediately follow it.

$anchorIndex = -1
$paraCount = $d.Paragraphs.Count
for ($i = 1; $i -le $paraCount; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*Người phỏng vấn thứ 5*") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -eq -1) {
    throw "Could not find the 'Người phỏng vấn thứ 5' heading paragraph."
}

# Họ & tên: Nguyễn Văn A -> Họ & tên: Lê Quang Tường
$rng = $d.Paragraphs.Item($anchorIndex + 1).Range
$rng.Find.Execute("Nguyễn Văn A", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "Lê Quang Tường", 2)

# Công việc: công việc -> Công việc: Kinh doanh
$rng = $d.Paragraphs.Item($anchorIndex + 2).Range
$rng.Find.Execute("công việc", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "Kinh doanh", 2)

# Tuổi: 000 -> Tuổi: 25
$rng = $d.Paragraphs.Item($anchorIndex + 3).Range
$rng.Find.Execute("000", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "25", 2)

# Giới tính: Nam, Nữ, Khác -> Giới tính: Nam
$rng = $d.Paragraphs.Item($anchorIndex + 4).Range
$rng.Find.Execute("Nam, Nữ, Khác", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "Nam", 2)
